# accounting-cover.docx: add forgotten applicant postcode/city/land lines
# and the Zahlungszweck line; flip the Normal style's overflowPunct on.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert "{{applicant_postcode}} {{applicant_city}}" paragraph right
#    after the existing "{{applicant_address}}" paragraph.
# ---------------------------------------------------------------------
$addressPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*{{applicant_address}}*") {
        $addressPara = $p
        break
    }
}

$addressPara.Range.InsertParagraphAfter()
$postcodePara = $addressPara.Next()
$postcodePara.Range.Text = "{{applicant_postcode}} {{applicant_city}}"

# ---------------------------------------------------------------------
# 2) Insert "{{applicant_land}}" paragraph right after that one.
# ---------------------------------------------------------------------
$postcodePara.Range.InsertParagraphAfter()
$landPara = $postcodePara.Next()
$landPara.Range.Text = "{{applicant_land}}"

# ---------------------------------------------------------------------
# 3) Insert the "Zahlungszweck:" line right after the "Fibukonto:" line.
# ---------------------------------------------------------------------
$fibuPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Fibukonto:*") {
        $fibuPara = $p
        break
    }
}

$fibuPara.Range.InsertParagraphAfter()
$zahlungszweckPara = $fibuPara.Next()
$zahlungszweckPara.Range.Text = "Zahlungszweck:" + [char]9 + [char]9 + "{{zahlungszweck}}"

# ---------------------------------------------------------------------
# 4) Flip the "Normal" style's overflowPunct (HangingPunctuation in the
#    Word object model) from false to true.
# ---------------------------------------------------------------------
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $true

Write-Output "edit complete"
